$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51 (shifts existing rows 51-87 down to 52-88)
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly record
$ws.Range("A51").Value = 3
$ws.Range("B51").Value = "Femacal de La Calera"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 44582
$ws.Range("D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E51").Value = 5
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100107
$ws.Range("H51").Value = "Otros"
$ws.Range("I51").Value = 100107011
$ws.Range("J51").Value = "Tuna"
$ws.Range("K51").Value = "Sin especificar"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 60
$ws.Range("N51").Value = 20000
$ws.Range("O51").Value = 20000
$ws.Range("P51").Value = 20000
$ws.Range("Q51").Value = "`$/caja 20 kilos"
$ws.Range("R51").Value = "Provincia de Limarí"
$ws.Range("S51").Value = 1000
$ws.Range("T51").Value = 20
